$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value2 = 131257424
$ws.Range("A7").Value2 = 131256691
$ws.Range("AB5").Value2 = "11:33"
$ws.Range("AB7").Value2 = "10:55"
$ws.Range("AC5").Value2 = "Gran"
$ws.Range("AC7").Value2 = "Ringhack på gran."
$ws.Range("B5").Value2 = 79245
$ws.Range("B7").Value2 = 57884
$ws.Range("E5").Value2 = 6425
$ws.Range("E7").Value2 = 100109
$ws.Range("F5").Value2 = "Garnlav"
$ws.Range("F7").Value2 = "Tretåig hackspett"
$ws.Range("G5").Value2 = "Alectoria sarmentosa"
$ws.Range("G7").Value2 = "Picoides tridactylus"
$ws.Range("H5").Value2 = "(Ach.) Ach."
$ws.Range("H7").Value2 = "(Linnaeus, 1758)"
$ws.Range("M5").Value2 = $null
$ws.Range("M7").Value2 = "äldre spår"
$ws.Range("Q5").Value2 = 488876
$ws.Range("Q7").Value2 = 488667
$ws.Range("R5").Value2 = 6665177
$ws.Range("R7").Value2 = 6665262
$ws.Range("Z5").Value2 = "11:33"
$ws.Range("Z7").Value2 = "10:55"
$ws.Range("A6").Value2 = 131255793
$ws.Range("A8").Value2 = 131260583
$ws.Range("AB6").Value2 = "09:56"
$ws.Range("AB8").Value2 = "15:30"
$ws.Range("AC6").Value2 = "Flera fruktkroppar."
$ws.Range("AC8").Value2 = "Ringhack på tall."
$ws.Range("B6").Value2 = 91830
$ws.Range("B8").Value2 = 57884
$ws.Range("E6").Value2 = 5432
$ws.Range("E8").Value2 = 100109
$ws.Range("F6").Value2 = "Granticka"
$ws.Range("F8").Value2 = "Tretåig hackspett"
$ws.Range("G6").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("G8").Value2 = "Picoides tridactylus"
$ws.Range("H6").Value2 = ""
$ws.Range("H8").Value2 = "(Linnaeus, 1758)"
$ws.Range("M6").Value2 = $null
$ws.Range("M8").Value2 = "färska spår"
$ws.Range("Q6").Value2 = 488817
$ws.Range("Q8").Value2 = 488834
$ws.Range("R6").Value2 = 6665110
$ws.Range("R8").Value2 = 6665228
$ws.Range("Z6").Value2 = "09:56"
$ws.Range("Z8").Value2 = "15:30"
$ws.Range("A12").Value2 = 131257290
$ws.Range("A14").Value2 = 131257520
$ws.Range("AB12").Value2 = "11:26"
$ws.Range("AB14").Value2 = "11:41"
$ws.Range("AC12").Value2 = "Ringhack på tall."
$ws.Range("AC14").Value2 = "Gran"
$ws.Range("B12").Value2 = 57884
$ws.Range("B14").Value2 = 79245
$ws.Range("E12").Value2 = 100109
$ws.Range("E14").Value2 = 6425
$ws.Range("F12").Value2 = "Tretåig hackspett"
$ws.Range("F14").Value2 = "Garnlav"
$ws.Range("G12").Value2 = "Picoides tridactylus"
$ws.Range("G14").Value2 = "Alectoria sarmentosa"
$ws.Range("H12").Value2 = "(Linnaeus, 1758)"
$ws.Range("H14").Value2 = "(Ach.) Ach."
$ws.Range("M12").Value2 = "äldre spår"
$ws.Range("M14").Value2 = $null
$ws.Range("Q12").Value2 = 488842
$ws.Range("Q14").Value2 = 488939
$ws.Range("R12").Value2 = 6665224
$ws.Range("R14").Value2 = 6665149
$ws.Range("Z12").Value2 = "11:26"
$ws.Range("Z14").Value2 = "11:41"
$ws.Range("A13").Value2 = 131256673
$ws.Range("A15").Value2 = 131260641
$ws.Range("AB13").Value2 = "10:54"
$ws.Range("AB15").Value2 = "15:34"
$ws.Range("AC13").Value2 = "Ringhack på tall."
$ws.Range("AC15").Value2 = "Ringhack på gran."
$ws.Range("Q13").Value2 = 488652
$ws.Range("Q15").Value2 = 488859
$ws.Range("R13").Value2 = 6665282
$ws.Range("R15").Value2 = 6665292
$ws.Range("Z13").Value2 = "10:54"
$ws.Range("Z15").Value2 = "15:34"
$ws.Range("A36").Value2 = 131260531
$ws.Range("A37").Value2 = 131257385
$ws.Range("AB36").Value2 = "15:25"
$ws.Range("AB37").Value2 = "11:31"
$ws.Range("AC36").Value2 = "Gran"
$ws.Range("AC37").Value2 = "Lågstubbe."
$ws.Range("B36").Value2 = 79245
$ws.Range("B37").Value2 = 91830
$ws.Range("E36").Value2 = 6425
$ws.Range("E37").Value2 = 5432
$ws.Range("F36").Value2 = "Garnlav"
$ws.Range("F37").Value2 = "Granticka"
$ws.Range("G36").Value2 = "Alectoria sarmentosa"
$ws.Range("G37").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("H36").Value2 = "(Ach.) Ach."
$ws.Range("H37").Value2 = ""
$ws.Range("Q36").Value2 = 488786
$ws.Range("Q37").Value2 = 488876
$ws.Range("R36").Value2 = 6665188
$ws.Range("R37").Value2 = 6665194
$ws.Range("Z36").Value2 = "15:25"
$ws.Range("Z37").Value2 = "11:31"
$ws.Range("A43").Value2 = 131273946
$ws.Range("A44").Value2 = 131273991
$ws.Range("Q43").Value2 = 488774
$ws.Range("Q44").Value2 = 488928
$ws.Range("R43").Value2 = 6665353
$ws.Range("R44").Value2 = 6665146
